$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F ("想去人数" / want-to-go count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1338
$ws1.Range("F4").Value = 14668
$ws1.Range("F5").Value = 17909
$ws1.Range("F6").Value = 17909
$ws1.Range("F8").Value = 78
$ws1.Range("F18").Value = 162
$ws1.Range("F19").Value = 46
$ws1.Range("F23").Value = 67
$ws1.Range("F25").Value = 7395
$ws1.Range("F27").Value = 9
$ws1.Range("F28").Value = 45
$ws1.Range("F29").Value = 1189
$ws1.Range("F30").Value = 13
$ws1.Range("F32").Value = 73
$ws1.Range("F34").Value = 147
$ws1.Range("F35").Value = 150
$ws1.Range("F37").Value = 5143
$ws1.Range("F38").Value = 22
$ws1.Range("F39").Value = 34

# Sheet "全部类型" (all types) - same underlying records, shifted by the
# extra rows present only in this sheet (23 and 32), column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1338
$ws4.Range("F4").Value = 14668
$ws4.Range("F5").Value = 17909
$ws4.Range("F6").Value = 17909
$ws4.Range("F8").Value = 78
$ws4.Range("F18").Value = 162
$ws4.Range("F19").Value = 46
$ws4.Range("F24").Value = 67
$ws4.Range("F26").Value = 7395
$ws4.Range("F28").Value = 9
$ws4.Range("F29").Value = 45
$ws4.Range("F30").Value = 1189
$ws4.Range("F31").Value = 13
$ws4.Range("F34").Value = 73
$ws4.Range("F36").Value = 147
$ws4.Range("F37").Value = 150
$ws4.Range("F39").Value = 5143
$ws4.Range("F40").Value = 22
$ws4.Range("F41").Value = 34
